# poland_inflation_ratings.xlsx — "Add files via upload"
#
# The underlying dataset (monthly Polish CPI inflation readings) was
# refreshed: the last figure for 2021 (M26, December) was filled in, and
# three more years of data (2022, 2023 and the single data point available
# so far for 2024) were appended below the existing table. A footer/header
# note and a refreshed cell selection are also applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fill in the previously-blank December 2021 figure.
# ---------------------------------------------------------------------
$ws.Cells.Item(26, 13).Value = 8          # M26

# ---------------------------------------------------------------------
# 2. Append the 2022 row (row 27), reusing the formatting of the row
#    immediately above it (style carries the "Year" column centering and
#    the data columns' vertical-centered / wrap formatting).
# ---------------------------------------------------------------------
$ws.Range("A26:M26").Copy()
$ws.Range("A27:M27").PasteSpecial(-4122)   # xlPasteFormats

$row2022 = @(2022, 8.7, 8.1, 10.2, 11.4, 12.8, 14.2, 14.2, 14.8, 15.7, 16.4, 16.1, 15.3)
for ($col = 1; $col -le $row2022.Length; $col++) {
    $ws.Cells.Item(27, $col).Value = $row2022[$col - 1]
}

# ---------------------------------------------------------------------
# 3. Append the 2023 row (row 28), same formatting approach.
# ---------------------------------------------------------------------
$ws.Range("A27:M27").Copy()
$ws.Range("A28:M28").PasteSpecial(-4122)

$row2023 = @(2023, 15.9, 17.2, 15.2, 14, 12.5, 11, 10.3, 9.5, 7.7, 6.3, 6.3, 6.2)
for ($col = 1; $col -le $row2023.Length; $col++) {
    $ws.Cells.Item(28, $col).Value = $row2023[$col - 1]
}

# ---------------------------------------------------------------------
# 4. Append the 2024 row (row 29) — only January is available so far.
# ---------------------------------------------------------------------
$ws.Range("A28:B28").Copy()
$ws.Range("A29:B29").PasteSpecial(-4122)

$ws.Cells.Item(29, 1).Value = 2024
$ws.Cells.Item(29, 2).Value = 3.9

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Confidentiality footer added to the sheet header.
# ---------------------------------------------------------------------
$ws.PageSetup.RightHeader = "&""Calibri""&10&K000000 Internal&1#`r"

# ---------------------------------------------------------------------
# 6. Refresh the active selection / scroll position to reflect the newly
#    entered data at the bottom of the table.
# ---------------------------------------------------------------------
[void]$ws.Activate()
try { [void]$ws.Cells.Item(7, 1).Select() } catch {}
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
[void]$ws.Range("B30").Select()

# ---------------------------------------------------------------------
# 7. Misc cosmetic/session settings (best-effort; harmless if the host
#    does not expose a given property).
# ---------------------------------------------------------------------
try { $ws.StandardHeight = 14.4 } catch {}
try { $wb.Styles.Item("Normalny").Name = "Normal" } catch {}
